$wb = $excel.ActiveWorkbook

# --- 1. Update status text "Ready for handoff" -> "In Translation" everywhere it appears ---
$oldText = "Ready for handoff"
$newText = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $cellText = [string]$cell.Text
            if ($oldText -eq $cellText) {
                $cell.Value = $newText
            }
        }
    }
}

# --- 2. Shrink the previously wide "status" columns ---
# New target column width (stored sheet width ~13.41 characters).
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
